# Weekly update: insert two new weekly price observations for
# "Terminal La Palmera de La Serena - Jengibre", pushing later rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 8 (old rows 8..21 shift down to 9..22) ---
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44428
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100114007
$ws.Cells.Item(8, 7).Value = "Jengibre"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 480
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14500
$ws.Cells.Item(8, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 1115
$ws.Cells.Item(8, 17).Value = 13
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# --- Insert second new row at row 16 (rows currently at 16.. shift down to 17..) ---
$ws.Rows.Item(16).Insert()

$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44435
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = "Jengibre"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 480
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 14000
$ws.Cells.Item(16, 13).Value = 13500
$ws.Cells.Item(16, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1038
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Dimension should now cover A1:R23
$ws.Range("A1:R23").Select()
